{"js": "// Replace the 25 \"two-digit \u00f7 one-digit\" answer strings that live in the\n// document's single practice table with a new set of problems/answers.\n//\n// The cells are processed in document order (row-major, left-to-right) and\n// matched positionally against `replacements` below \u2014 this is important\n// because a couple of the \"before\" strings are not unique (e.g. \"66\u00f73=22, 0\"\n// appears twice) but map to two *different* \"after\" strings depending on\n// which occurrence it is. Matching by position (rather than a global\n// find/replace-all on text) keeps each occurrence mapped to the right\n// target. Blank spacer cells/rows in the table are skipped.\nconst replacements = [\n  [\"70\u00f73=23, 1\", \"63\u00f73=21, 0\"],\n  [\"66\u00f73=22, 0\", \"77\u00f79=8, 5\"],\n  [\"48\u00f72=24, 0\", \"61\u00f74=15, 1\"],\n  [\"96\u00f77=13, 5\", \"12\u00f73=4, 0\"],\n  [\"39\u00f77=5, 4\", \"56\u00f79=6, 2\"],\n  [\"80\u00f77=11, 3\", \"55\u00f74=13, 3\"],\n  [\"86\u00f77=12, 2\", \"95\u00f76=15, 5\"],\n  [\"43\u00f78=5, 3\", \"40\u00f74=10, 0\"],\n  [\"64\u00f76=10, 4\", \"39\u00f72=19, 1\"],\n  [\"32\u00f73=10, 2\", \"24\u00f75=4, 4\"],\n  [\"46\u00f78=5, 6\", \"52\u00f74=13, 0\"],\n  [\"40\u00f74=10, 0\", \"56\u00f76=9, 2\"],\n  [\"64\u00f74=16, 0\", \"59\u00f78=7, 3\"],\n  [\"19\u00f75=3, 4\", \"37\u00f76=6, 1\"],\n  [\"66\u00f73=22, 0\", \"28\u00f74=7, 0\"],\n  [\"38\u00f74=9, 2\", \"54\u00f75=10, 4\"],\n  [\"82\u00f73=27, 1\", \"57\u00f79=6, 3\"],\n  [\"91\u00f78=11, 3\", \"68\u00f76=11, 2\"],\n  [\"52\u00f75=10, 2\", \"15\u00f75=3, 0\"],\n  [\"46\u00f74=11, 2\", \"72\u00f79=8, 0\"],\n  [\"35\u00f78=4, 3\", \"28\u00f77=4, 0\"],\n  [\"45\u00f73=15, 0\", \"40\u00f76=6, 4\"],\n  [\"85\u00f73=28, 1\", \"32\u00f75=6, 2\"],\n  [\"96\u00f74=24, 0\", \"67\u00f74=16, 3\"],\n  [\"57\u00f76=9, 3\", \"78\u00f72=39, 0\"],\n];\n\n// The document has exactly one table; walk it row by row, cell by cell.\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.load(\"text\");\n  }\n}\nawait context.sync();\n\nlet i = 0;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    const text = cell.body.text;\n    // Spacer rows contain empty cells \u2014 nothing to replace there.\n    if (text.length === 0) {\n      continue;\n    }\n    if (i >= replacements.length) {\n      break;\n    }\n    const [oldText, newText] = replacements[i];\n    if (text !== oldText) {\n      throw new Error(\n        \"Mismatch at occurrence \" + i + \": expected '\" + oldText + \"' but found '\" + text + \"'\"\n      );\n    }\n    // Replace the text of the cell's first (only) paragraph in place so the\n    // existing run/paragraph formatting (font, size, alignment) is kept \u2014\n    // setting cell.body text directly would drop the run's rPr/pPr.\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n    const range = paragraphs.items[0].getRange();\n    range.insertText(newText, \"Replace\");\n    i++;\n  }\n}\nawait context.sync();\n\nif (i !== replacements.length) {\n  throw new Error(\n    \"Expected to replace \" + replacements.length + \" cells but only matched \" + i\n  );\n}\n", "ps1": "# Replace the 25 \"two-digit \u00f7 one-digit\" answer strings that live in the\n# document's single practice table with a new set of problems/answers.\n#\n# Cells are processed in document order (row-major, left-to-right) and\n# matched positionally against $replacements below -- this matters because a\n# couple of the \"before\" strings are not unique (e.g. \"66\u00f73=22, 0\" appears\n# twice) but map to two *different* \"after\" strings depending on which\n# occurrence it is. Matching by position (instead of Find/Replace on text)\n# keeps each occurrence mapped to the correct target. Blank spacer\n# rows/cells in the table are skipped.\n$replacements = @(\n    @('70\u00f73=23, 1', '63\u00f73=21, 0'),\n    @('66\u00f73=22, 0', '77\u00f79=8, 5'),\n    @('48\u00f72=24, 0', '61\u00f74=15, 1'),\n    @('96\u00f77=13, 5', '12\u00f73=4, 0'),\n    @('39\u00f77=5, 4', '56\u00f79=6, 2'),\n    @('80\u00f77=11, 3', '55\u00f74=13, 3'),\n    @('86\u00f77=12, 2', '95\u00f76=15, 5'),\n    @('43\u00f78=5, 3', '40\u00f74=10, 0'),\n    @('64\u00f76=10, 4', '39\u00f72=19, 1'),\n    @('32\u00f73=10, 2', '24\u00f75=4, 4'),\n    @('46\u00f78=5, 6', '52\u00f74=13, 0'),\n    @('40\u00f74=10, 0', '56\u00f76=9, 2'),\n    @('64\u00f74=16, 0', '59\u00f78=7, 3'),\n    @('19\u00f75=3, 4', '37\u00f76=6, 1'),\n    @('66\u00f73=22, 0', '28\u00f74=7, 0'),\n    @('38\u00f74=9, 2', '54\u00f75=10, 4'),\n    @('82\u00f73=27, 1', '57\u00f79=6, 3'),\n    @('91\u00f78=11, 3', '68\u00f76=11, 2'),\n    @('52\u00f75=10, 2', '15\u00f75=3, 0'),\n    @('46\u00f74=11, 2', '72\u00f79=8, 0'),\n    @('35\u00f78=4, 3', '28\u00f77=4, 0'),\n    @('45\u00f73=15, 0', '40\u00f76=6, 4'),\n    @('85\u00f73=28, 1', '32\u00f75=6, 2'),\n    @('96\u00f74=24, 0', '67\u00f74=16, 3'),\n    @('57\u00f76=9, 3', '78\u00f72=39, 0')\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cellRange = $cell.Range\n        $fullText = $cellRange.Text\n        # A cell's Range.Text always ends with the cell-mark (and, because\n        # each cell here holds a single paragraph, a paragraph mark before\n        # it) -- strip those last two control characters to get the content.\n        $text = $fullText.Substring(0, $fullText.Length - 2)\n        if ($text.Length -eq 0) {\n            # Blank spacer cell -- nothing to do.\n            continue\n        }\n        if ($i -ge $replacements.Length) {\n            break\n        }\n        $pair = $replacements[$i]\n        $old = $pair[0]\n        $new = $pair[1]\n        if ($text -ne $old) {\n            throw \"Mismatch at row $r col $c : expected '$old' but found '$text'\"\n        }\n        # Assign only to the cell's Range (not the whole-cell End-of-cell\n        # mark) so Word replaces the run's text while keeping its existing\n        # character/paragraph formatting (font, size, alignment).\n        $cellRange.Text = $new\n        $i++\n    }\n}\n\nif ($i -ne $replacements.Length) {\n    throw \"Expected to replace $($replacements.Length) cells but only replaced $i\"\n}\n"}
